$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

# "Figures are for 16+ year olds." -> "Figures are for 16-64 year olds."
# (shared footnote text used by several rows in column F)
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*Figures are for 16+ year olds.*") {
        $cell.Value2 = $val -replace "Figures are for 16\+ year olds\.", "Figures are for 16-64 year olds."
    }
}

# F6 footnote: "Employment volumes are for 16+ year olds." -> "...16-64 year olds."
# and add a new bullet about industry/occupation split volumes right after it.
$f6 = $ws.Cells.Item(6, 6)
$f6val = $f6.Value2
$f6val = $f6val -replace "Employment volumes are for 16\+ year olds\.", "Employment volumes are for 16-64 year olds."
$f6val = $f6val -replace "(Employment volumes are for 16-64 year olds\.</li>)\r?\n", "`$1`n<li>Industry and occupation split volumes are for all ages.</li>`n"
$f6.Value2 = $f6val

# That row's text now wraps across more lines, so grow the row to fit.
$ws.Rows.Item(6).RowHeight = 276

# Scroll the view back to the top of the sheet and reselect F7.
$ws.Activate()
$ws.Range("F7").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
